# special_case.xlsx - "modify for interval & last"
# Append a third data row (A006 / TASK1 / UT) below the existing two rows,
# then move the active selection down to A4, mirroring what a user does
# after typing a new row by hand and pressing Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "A006"
$ws.Range("B3").Value = "TASK1"
$ws.Range("C3").Value = "UT"

$ws.Range("A4").Select()
